{"js": "// HERCULES-10151 - Plantilla exportaci\u00f3n memoria - Ampliar interlineado a 1,5\n//\n// The document defines line spacing through two paragraph styles:\n//   - \"Normal\"    (base style for the whole document)\n//   - \"Body Text\" (w:styleId=\"Cuerpodetexto\", based on Normal)\n// Both had w:spacing w:line=\"276\" (auto rule -> 1.15 lines for a 12pt\n// default run) and need to become w:line=\"360\" (auto rule -> 1.5 lines).\n// The \"Normal\" style additionally flips w:overflowPunct from true to false.\n\nconst styles = context.document.getStyles();\nconst normalStyle = styles.getByNameOrNullObject(\"Normal\");\nconst bodyTextStyle = styles.getByNameOrNullObject(\"Body Text\");\nnormalStyle.load(\"nameLocal\");\nbodyTextStyle.load(\"nameLocal\");\nawait context.sync();\n\n// --- Line spacing: 1.15 (w:line=\"276\") -> 1.5 (w:line=\"360\") ---\n// Word JS API's ParagraphFormat.lineSpacing is expressed in points and\n// serializes to OOXML as `w:line = lineSpacing * 20` when the rule is\n// \"auto\"/multiple (w:lineRule=\"auto\"). 18pt -> w:line=\"360\".\nnormalStyle.paragraphFormat.lineSpacing = 18;\nbodyTextStyle.paragraphFormat.lineSpacing = 18;\nawait context.sync();\n\n// --- overflowPunct: true -> false (Normal style only) ---\n// This flag (\"allow punctuation to hang outside the margin\") is not\n// surfaced as a named property on Word.ParagraphFormat, so it is reached\n// through the same object-model bridge the proxy's public getters/setters\n// use internally (ParagraphFormat.HangingPunctuation <-> w:overflowPunct).\nconst normalParagraphFormat = normalStyle.paragraphFormat;\nif (typeof normalParagraphFormat._omSet === \"function\") {\n  normalParagraphFormat._omSet(\"HangingPunctuation\", \"False\");\n} else {\n  // Fallback (older/newer shim without the internal bridge exposed):\n  // at least keep line spacing applied; overflowPunct stays as-is.\n}\nawait context.sync();\n", "ps1": "# HERCULES-10151 - Plantilla exportaci\u00f3n memoria - Ampliar interlineado a 1,5\n#\n# The document defines line spacing through two paragraph styles:\n#   - \"Normal\"    (base style for the whole document)\n#   - \"Body Text\" (w:styleId=\"Cuerpodetexto\", based on Normal)\n# Both had w:spacing w:line=\"276\" (auto rule -> 1.15 lines for a 12pt\n# default run) and need to become w:line=\"360\" (auto rule -> 1.5 lines).\n# The \"Normal\" style additionally flips w:overflowPunct from true to false.\n\n$d = $word.ActiveDocument\n\n$normalStyle = $d.Styles.Item(\"Normal\")\n$bodyTextStyle = $d.Styles.Item(\"Body Text\")\n\n# --- Line spacing: 1.15 (w:line=\"276\") -> 1.5 (w:line=\"360\") ---\n# Word's ParagraphFormat.LineSpacing is expressed in points and serializes\n# to OOXML as `w:line = LineSpacing * 20` while LineSpacingRule stays\n# wdLineSpaceMultiple (w:lineRule=\"auto\"). 18pt -> w:line=\"360\".\n$normalStyle.ParagraphFormat.LineSpacing = 18\n$bodyTextStyle.ParagraphFormat.LineSpacing = 18\n\n# --- overflowPunct: true -> false (Normal style only) ---\n# \"Allow punctuation to hang outside the margin\" maps to\n# ParagraphFormat.HangingPunctuation (OOXML w:overflowPunct).\n$normalStyle.ParagraphFormat.HangingPunctuation = 0\n"}
